# Sheet1 holds the wealth-tax model: h (B2), m (B3), f (B4) are the
# calibrated model parameters/scale factor; everything else on the sheet
# (columns D:G, the simulation block in rows 36-37/42-43, and the three
# scatter charts) derives from these via formulas and recalculates on its
# own once the inputs change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calibrated wealth-tax parameters (h, m) and the wealth scale factor (f)
$ws.Range("B2").Value = 1.1867121650630199
$ws.Range("B3").Value = 1.12835012204532
$ws.Range("B4").Value = 210416.55648500001

# Simulation conversion anchor (named range sim_conv -> Sheet1!$A$42),
# kept in sync with f by hand since it is a plain literal, not a formula.
$ws.Range("A42").Value = 210416.55648500001

# Finished baseline simulation: updated wealth inputs feeding row 41
$ws.Range("B41").Value = 0.63388029999999995
$ws.Range("C41").Value = 0.77039457
$ws.Range("D41").Value = 0.87985102999999998
$ws.Range("E41").Value = 0.98640106999999999
$ws.Range("F41").Value = 1.10487501
$ws.Range("G41").Value = 1.52812568
$ws.Range("H41").Value = 2.66742334

# Leave the selection where the author left it after finishing up
$ws.Range("B3").Select()
